# Apply the updated crypto price/volume figures scraped on
# Sun Sep 29 22:36:45 UTC 2024 (GitHub Actions refresh of cryptos.xlsx).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.777.73'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').Value = '2.663.95'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '599.35'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '159.97'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.01%  '
$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.642'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +4.02%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  -2.20%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '5.86'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.48%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.399'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.15%  '
$ws.Range('E12').Value = '  +1.32%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '29.18'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.48%  '
$ws.Range('E14').Value = '  -1.42%  '
$ws.Range('D15').Value = '3.141.25'
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('D16').Value = '65.681.94'
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('D17').Value = '2.682.65'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '12.57'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.72%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.79'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '354.55'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.63%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.46'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.53%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.999'
$ws.Range('D22').ClearFormats()
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '69.88'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('E24').Value = '  +8.58%  '
$ws.Range('E25').Value = '  +0.19%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.77'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.22%  '
$ws.Range('E27').Value = '  +1.97%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '565.90'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +6.17%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.13'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.08%  '
$ws.Range('E30').Value = '  -2.24%  '
$ws.Range('E31').Value = '  +0.12%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.15'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.53%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.83'
$ws.Range('D33').ClearFormats()
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.71'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +3.38%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.51'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.36%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.423'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '20.59'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('E38').Value = '  +1.87%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '154.32'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -3.19%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.51'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +7.34%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '161.87'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.15%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '4.10'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.91%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0617'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.05%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '23.51'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.19%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.645'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.00%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0259'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.102'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.72%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '19.78'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.23%  '
$ws.Range('E50').Value = '  -7.31%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.815'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.40%  '
